# "editando labels de imagens"
#
# The followers/likes figures for the image014-image029 block (rows 13-30)
# were re-shuffled: rows 13-15 worth of (followers, likes) numbers were
# removed while the "imageNNN" labels in column A stayed put, and the last
# three rows of the table (28-30, labels image027-image029) were deleted
# outright. That both shrinks the used range down to A1:C27 and drops the
# now-orphaned "image027".."image029" entries out of the shared-string
# table on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New followers (B) / likes (C) values for rows 13-27; column A keeps its
# original "imageNNN" label in every one of these rows.
$newValues = @{
    13 = @(2800000, 107940)
    14 = @(4100000, 28825)
    15 = @(4100000, 29575)
    16 = @(4100000, 47067)
    17 = @(4100000, 65789)
    18 = @(102000,  2585)
    19 = @(102000,  3438)
    20 = @(102000,  3955)
    21 = @(989000,  7593)
    22 = @(989000,  8279)
    23 = @(989000,  10571)
    24 = @(989000,  5307)
    25 = @(989000,  6500)
    26 = @(989000,  7593)
    27 = @(989000,  9142)
}

foreach ($row in $newValues.Keys) {
    $pair = $newValues[$row]
    $ws.Range("B$row").Value = $pair[0]
    $ws.Range("C$row").Value = $pair[1]
}

# Drop the trailing rows that no longer belong to the table.
$ws.Rows("28:30").Delete()

# Match the saved selection/scroll state (active cell A12, scrolled so row 4
# is at the top).
$ws.Range("A12").Select()
$excel.ActiveWindow.ScrollRow = 4
